$wb = $excel.ActiveWorkbook

# --- Append Week 15 simulation data points to the long-running simulation result lists ---
$ws = $wb.Worksheets.Item("YDS")
$s_YDS_B2 = @"
5 8 5 8 1 1 3 5 4 2 0 21 12 9 2 -1 16 14 4 9 2 2 3 8 5 0 3 16 4 1 1 4 -2 11 8 1 7 8 3 1 6 8 3 2 9 3 11 4 2 39 3 1 4 4 1 3 3 -2 12 5 4 39 -3 2 0 13 1 14 6 1 10 15 1 7 12 0 11 1 2 1 6 2 2 5 2 2 8 1 6 4 3 0 4 2 7 7 3 4 2 7 13 2 6 5 3 1 15 4 3 9 2 0 4 -2 1 3 10 8 2 10 10 3 0 0 8 3 4 2 0 6 7 6 3 16 2 3 2 3 20 7 3 4 12 -2 2 1 25 6 4 4 12 3 7 -3 1 -1 2 16 1 3 3 -4 4 4 2 12 0 2 7 21 3 4 9 7 3 1 37 2 4 2 3 0 1 0 1 3 17 3 3 2 7 -2 6 4 -1 2 5 13 5 5 9 11 5 0 7 1 5 13 14 4 5 12 4 2 0 1 4 8 70 5 9 4 4 13 3 14 18 1 1 5 4 0 4 5 2 0 3 2 1 5 0 3 -1 6 2 0 0 2 3 1 0 1 8 4 13 6 0 7 0 14 4 1 3 1 4 4 2 11 3 5 3 1 1 -1 2 2 9 2 1 4 1 -1 5 13 5 5 8 7 6 16 1 6 0 0 4 6 2 9 0 6 2 6 -1 2 5 11 5 2 4 4 10 1 14 3 0 7 6 2 -3 1 4 4 10 0 2 4 7 1 0 10 2 1 3 12 6 10 0 3 0 1 1 7 3 4 6 5 -1 10 4 3 6 3 1 5 13 7 2 5 3 8 2 5 1 3 4 9 4 4 1 12 3 14 0 3 4 0 4 0 4 4 10 7 1 12 7 0 10 0 13 7 0 6 6 2 7 20 4 10 11 8 4 18 1 11 2 1 2 2 9 1 3 0 15 2 15 13 7 2 8 1 3 4 -2 0 5 11 7 2 -3 5 2 1 6 4 1 5 5 20 1 2 5 2 4 1 4 7 4 19 5 4 1 0 0 1 5 -1 1 3 1 2 17 2 7 -3 4 1 13 0 8 -1 6 5 2 1 0 -1 0 0 11 2 16 14 6 7 3 0 15 2 0 29 9 6 12 0 10 0 13 2 3 0 -2 8 1 1 9 7 1 3 1 4 8 2 2 12 2 2 9 3 4 3 6 3 -1 5 0 2 24 9 2 3 4 13 2 0 1 1 -2 2 3 1 11 2 4 1 0 1 10 5 -1 -2 2 1 6 6 -2 13 2 2 2 4 -2 -2 9 4 0 5 5 11 3 6 0 3 3 3 0 5 48 2 4 -1 3 2 2 1 2 1 0 -3 2 -3 3 4 14 5 9 1 3 4 3 5 1 8 9 6 0 14 2 0 30 16 -1 -1 23 4 0 3 3 -1 -3 16 5 4 16 6 2 -1 15 0 2 3 2 0 0 13 19 4 3 -2 8 4 8 8 7 1 4 -4 11 15 -3 4 66 1 3 3 0 0 1 2 -1 1 0 -1 0 -8 1 24 1 3 -1 5 3 -4 12 6 6 8 5 5 3 4 7 4 -3 5 -1 3 1 4 0 2 7 0 1 5 2 4 1 11 -1 6 -1 3 3 1 2 5 0 -1 4 0 7 4 3 5 13 2 4 0 7 2 0 4 3 7 12 4 1 3 1 1 0 4 7 4 30 0 6 -4 7 1 3 1 2 7 5 1 4 9 3 5 -4 3 10 -1 2 5 1 8 11 4 8 0 1 3 4 4 2 5 1 20 4 30 14 3 3 16 15 3 29 9 3 0 1 7 7 9 8 5 -4 7 14 -3 5 -2 8 5 4 1 11 2 4 -1 0 2
"@
$ws.Range("B2").Value = $s_YDS_B2
$s_YDS_B3 = @"
7 25 18 12 22 9 -2 5 29 37 10 11 17 22 6 7 9 9 5 12 19 0 18 13 3 3 22 2 6 24 17 5 11 9 12 8 6 3 31 16 33 15 71 7 11 3 7 8 18 29 11 3 8 26 8 5 5 39 23 19 9 25 25 7 23 11 5 13 7 7 9 9 7 10 8 14 9 3 -5 3 25 9 10 14 3 12 11 11 6 10 -2 4 4 36 19 18 14 9 35 19 8 6 11 11 17 1 2 35 9 2 19 12 49 9 12 13 16 7 12 6 16 12 50 7 22 22 30 9 17 8 10 29 22 8 35 1 7 2 8 9 12 12 6 17 23 -3 4 9 9 13 11 14 4 2 8 54 7 22 14 6 8 21 10 28 12 11 8 27 4 11 26 7 51 6 2 8 13 20 7 3 13 39 12 4 5 12 6 5 2 12 3 8 -3 6 9 7 12 6 5 9 17 9 16 0 7 15 18 6 17 5 6 17 10 15 6 7 12 25 10 10 5 7 14 11 12 20 3 5 1 5 10 11 11 40 12 9 12 20 5 14 4 11 5 18 9 12 9 10 3 15 15 2 6 17 4 5 15 15 14 8 9 5 5 14 25 9 2 4 5 8 21 23 3 9 5 9 2 3 16 26 30 13 3 9 14 5 9 20 6 12 7 11 7 16 11 9 7 19 13 25 8 13 8 14 5 1 2 4 6 24 9 4 5 14 21 11 19 14 14 6 13 15 9 25 11 19 14 4 20 9 8 4 -3 14 10 40 15 38 15 7 14 30 13 8 5 9 6 12 11 9 -1 2 9 4 17 11 25 5 5 6 20 0 9 10 10 5 6 9 8 34 10 7 27 0 6 14 7 3 5 12 64 6 14 5 3 7 9 14 -3 4 15 6 6 12 13 7 13 3 17 12 9 8 20 17 8 10 7 6 23 26 7 15 7 4 7 9 16 3 -4 15 5 8 9 -1 21 11 28 4 14 15 2 11 9 1 11 22 4 12 20 2 3 1 6 20 17 14 -3 31 5 3 19 6 37 14 4 17 6 3 6 14 22 2 15 4 9 2 11 6 3 14 9 7 16 8 21 6 19 11 9 6 7 7 3 3 20 7 24 13 2 -6 8 14 7 9 24 35 5 17 26 40 3 4 12 6 8 6 7 13 -3 27 31 3 20 3 5 7 -2 6 11 18 7 7 -4 -2 4 4 9 -3 6 13 6 3 6 50 5 3 7 7 20 4 18 5 9 5 6 27 1 2 14 6 -1 10 24 7 7 9 3 18 21 21 5 4 27 19 3 12 15 16 3 27 5 18 5 43 3 18 56 15 10 5 7 15 6 7 9 21 5 19 6 7 7 3 23 6 19 26 2 11 29 7 4 30 15 20 7 6 14 15 24 16 5 17 7 5 23 2 19 0 18 -2 -2 34 6 7 14 6 2 16 15 11 5 19 10 48 7 16 4 21 14 16 5 14 5 -2 11 3 8 5 8 14 18 22 19 17 13 0 1 62 12 17
"@
$ws.Range("B3").Value = $s_YDS_B3
$s_YDS_C2 = @"
2 4 9 2 5 19 5 2 0 6 4 0 1 1 4 15 3 3 2 1 8 21 3 5 7 12 5 7 13 6 -2 5 3 3 7 4 1 2 5 5 7 4 3 1 5 1 3 -1 3 9 4 7 3 9 5 -1 4 5 0 5 9 3 1 5 3 10 0 0 7 4 2 4 0 9 1 -3 16 2 -3 0 7 12 5 6 5 2 1 3 1 5 12 5 3 5 5 2 4 -1 5 2 6 5 4 4 1 8 5 3 15 4 1 2 1 7 2 6 6 1 3 0 0 5 1 19 5 4 7 12 3 6 0 6 29 6 3 1 1 0 17 9 1 -4 6 8 -1 5 4 -2 0 8 4 9 -3 -4 2 -2 2 4 4 7 4 8 -1 2 4 3 6 -2 0 3 9 0 0 4 4 6 10 9 4 3 3 7 2 3 8 1 2 6 4 8 3 5 1 2 9 1 7 8 4 4 2 8 16 5 5 4 9 1 2 2 -4 2 5 2 2 17 6 -1 2 4 4 1 5 14 5 5 3 2 4 6 2 6 0 0 0 -5 5 -1 -2 3 2 2 5 3 4 7 5 9 4 6 7 4 5 0 1 4 3 1 3 15 6 42 4 5 2 18 3 6 2 10 6 5 13 10 3 13 6 2 3 -4 3 -2 4 3 -1 1 -2 8 3 4 1 9 5 -4 4 2 9 6 4 4 9 1 13 4 3 3 17 4 5 6 2 2 3 -1 2 4 3 3 1 15 2 5 2 5 3 1 5 3 3 3 1 1 10 5 3 14 4 2 16 3 -3 7 14 2 19 3 8 2 16 8 10 9 0 1 3 7 3 2 1 2 6 3 4 2 7 6 -1 3 11 14 2 8 3 5 1 1 7 3 9 3 2 -1 4 0 1 11 3 12 40 11 7 6 1 1 21 9 11 4 7 5 4 9 1 3 3 5 9 4 2 3 5 1 17 6 4 7 2 3 7 2 5 2 4 1 4 3 4 3 4 5 3 38 4 1 6 4 6 -1 5 4 6 3 2 13 8 5 2 5 -2 3 4 1 2 -1 13 1 3 6 4 7 2 2 5 10 19 1 5 5 9 12 -1 7 1 2 8 2 4 0 3 0 5 0 14 0 9 6 -1 6 2 8 12 8 9 3 4 6 3 1 6 6 -3 4 2 7 4 2 2 1 8 8 1 11 8 3 3 30 5 6 1 4 4 3 4 5 3 4 2 4 5 2 3 2 5 2 9 6 1 0 1 33 3 12 13 2 1 6 8 1 2 6 2 2 9 1 13 0 2 3 2 3 12 1 1 11 4 9 5 4 3 1 2 4 13 4 9 0 3 3 0 4 3 7 3 2 8 5 3 6 2 3 5 -2 6 9 6 0 5 2 5 2 6 30 0 13 8 2 6 2 3 5 4 2 7 2 5 2 2 2 6 2 4 2 2 4 2 1 2 1 5 8 5 6 4 4 2 6 3 4 9 3 5 6 7 11 6 12 0 1 7 1 4 9 4 6 1 3 8 7 8 5 10 1 6 5 4 2 3 19 5 12 1 5 4 5 1 0 4 3 17 1 1 3 5 5 9 2 6 3 3 5 4 1 5 3 2 1 0 11 2 5 4 18 11 7 6 5 7 4 4 3 -3 6 15 20 10 1 3 4 9 3 4 10 4 4 2 5 49 4 3 2 4 12 0 1 3 -6 5 1 0 6 7 9 4 3 1 2 10 1 2 8 4 5 4 6 8 2 4 3 0 2 -2 2 11 1 0 5 9 3 5 6 0 2 3 -2 7 2 -1 3 -2 9 3 4 4 23 8 3 11 13 1 1 2 3 6 -2
"@
$ws.Range("C2").Value = $s_YDS_C2
$s_YDS_C3 = @"
6 6 7 9 15 -1 -1 14 12 7 10 1 6 3 8 24 45 8 1 12 3 39 1 38 1 4 4 40 6 5 3 7 16 9 16 4 14 33 4 5 21 -7 7 14 27 2 28 6 44 -1 16 0 13 10 9 9 12 23 5 6 5 5 21 6 38 61 6 8 11 14 8 1 6 10 8 15 14 26 13 7 6 24 43 4 14 19 0 13 23 15 16 7 13 2 3 12 -6 9 13 20 19 2 13 8 1 17 7 17 15 3 20 23 8 5 17 3 3 18 14 9 4 5 8 3 17 4 3 17 7 27 7 6 7 14 31 8 5 8 5 4 1 5 16 14 6 5 1 45 19 6 6 7 15 5 5 2 26 27 7 7 7 17 1 5 -1 13 0 9 6 4 10 5 12 5 20 8 11 15 9 16 14 9 12 9 19 28 2 12 11 8 13 9 11 0 7 7 3 5 24 9 8 7 21 3 -1 7 10 -4 18 21 6 4 14 19 4 13 3 6 5 10 32 10 6 9 7 -3 10 5 10 2 4 17 18 1 5 41 7 25 11 6 10 5 34 9 13 35 4 12 24 8 28 6 7 34 1 4 7 6 6 13 6 3 4 5 6 6 5 6 19 4 7 4 10 9 22 21 9 20 48 16 6 4 3 9 7 29 2 11 10 6 7 4 24 8 14 12 10 5 16 9 16 13 22 7 35 9 9 19 8 19 19 31 7 11 19 1 9 26 3 44 23 7 6 4 4 41 -1 9 23 43 3 15 20 38 2 7 26 6 31 18 21 1 -1 16 6 10 5 3 7 12 17 13 2 14 2 9 50 4 28 5 14 15 18 5 6 4 8 15 16 4 0 15 3 34 9 19 77 4 18 6 13 25 9 11 29 6 2 6 17 2 3 7 17 5 28 10 10 12 9 5 21 7 39 7 12 11 31 5 14 6 -3 19 7 17 15 3 11 21 8 3 3 9 8 9 12 11 5 9 24 19 10 8 13 4 13 8 4 4 6 21 8 14 5 2 17 2 8 7 18 21 10 9 0 16 18 17 5 15 14 1 4 4 21 6 4 7 5 8 17 8 20 29 4 9 3 7 73 10 2 35 1 13 15 7 19 33 18 1 7 -3 15 5 13 13 15 7 -1 9 4 5 12 22 3 3 8 8 18 11 5 22 6 8 7 7 20 20 7 13 7 11 10 18 2 7 8 5 4 22 17 11 10 2 7 7 26 13 6 2 37 2 4 13 3 5 39 8 17 12 15 8 25 13 12 17 10 9 13 26 4 18 75 21 37 3 2 11 2 9 24 3 13 13 30 12 8 11 5 26 9 9 25 9 28 23 6 5 27 7 6 24 7 8 14 9 5 6 13 9 9 9 10 8 7 6 15 6 5 6 6 17 29 5 -2 8 3 14 1 37 30 6 15 6 38 9 1 13 3 10 12
"@
$ws.Range("C3").Value = $s_YDS_C3

$ws = $wb.Worksheets.Item("ST")
$s_ST_B4 = @"
54 54 67 66 62 64 67 63 64 62 51 46 39 69 43 59 69 56 44 57 66 64 63 49 40 68 69 66 45 62 50 65 69 67 61 67 64 69 65 67
"@
$ws.Range("B4").Value = $s_ST_B4
$s_ST_B5 = @"
18 18 24 30 29 22 19 29 27 26 11 11 2 4 7 31 16 9 18 20 30 33 24 17 15 70 24 11 14 21 21 26 24 22 24 14 30 23 33 17
"@
$ws.Range("B5").Value = $s_ST_B5
$s_ST_D3 = @"
40 40 36 46 43 43 53 33 46 37 42 46 41 39 57 49 53 46 50 35 35 43 52 51 54 46 40 48 40 41 59 54 48 51 54 50 49 46 45 40 49 52 28 48 43 37 51 34 46 55 49 50 38 39 46 49 61 52 41 51 63 46 43 50 60 32 49 45 58 47 58 46 36 46 50 54 46 44 54 53 37 44 37 37 39 51 29 53 54 41 53 41 40 42 55 45 33 47 45 56 59 46 49 37 46 49 51
"@
$ws.Range("D3").Value = $s_ST_D3
$s_ST_D4 = @"
0 0 11 12 16 11 3 0 2 0 5 0 13 0 15 0 24 0 12 0 0 9 21 0 12 0 12 32 0 18 10 20 0 11 20 13 4 0 0 0 0 1 0 4 1 0 0 0 1 16 10 12 0 8 15 0 12 15 3 0 0 6 0 8 5 0 0 0 8 8 9 0 2 -2 14 0 0 1 0 4 0 0 7 0 -3 0 0 0 0 0 14 0 9 14 13 0 0 0 0 12 2 0 0 0 0 8 0
"@
$ws.Range("D4").Value = $s_ST_D4
$s_ST_D5 = @"
0 0 0 0 0 0 0 6 0 0 0 0 9 0 0 0 0 0 -3 0 0 0 0 0 3 0 0 0 0 7 8 0 0 6 10 5 0 0 13 7 0 0 0 0 0 0 12 18 0 -4 15 13 -1 0 11 0 0 2 0 0 0 0 2 0 7 0 0 0 0 0 7 0 0 0 0 0 45 0 0 7 17 0 8 0 0 0 16 2 5 16 12 0 0 15 0
"@
$ws.Range("D5").Value = $s_ST_D5

# --- Update aggregate/summary numeric cells to reflect the newly added week ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 384
$ws.Range("D2").Value = 26
$ws.Range("F2").Value = 114
$ws.Range("G2").Value = 123
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 52
$ws.Range("L2").Value = 478
$ws.Range("M2").Value = 320
$ws.Range("O2").Value = 40
$ws.Range("P2").Value = 23
$ws.Range("Q2").Value = 897
$ws.Range("B3").Value = 21
$ws.Range("C3").Value = 302
$ws.Range("E3").Value = 59
$ws.Range("F3").Value = 163
$ws.Range("G3").Value = 58
$ws.Range("H3").Value = 60
$ws.Range("I3").Value = 104
$ws.Range("J3").Value = 83

$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 396
$ws.Range("D2").Value = 16
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 119
$ws.Range("G2").Value = 118
$ws.Range("J2").Value = 55
$ws.Range("L2").Value = 501
$ws.Range("M2").Value = 333
$ws.Range("O2").Value = 41
$ws.Range("P2").Value = 28
$ws.Range("Q2").Value = 932
$ws.Range("B3").Value = 22
$ws.Range("C3").Value = 302
$ws.Range("E3").Value = 40
$ws.Range("F3").Value = 176
$ws.Range("G3").Value = 82
$ws.Range("H3").Value = 40
$ws.Range("I3").Value = 98
$ws.Range("J3").Value = 87
$ws.Range("N3").Value = 29

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 158
$ws.Range("D2").Value = 107
$ws.Range("F2").Value = 69
$ws.Range("G2").Value = 61
$ws.Range("J2").Value = 65
$ws.Range("K2").Value = 54
$ws.Range("L2").Value = 35
$ws.Range("M2").Value = 22
$ws.Range("N2").Value = 14
$ws.Range("B3").Value = 118

$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("B2").Value = 14
$ws.Range("C2").Value = 11
$ws.Range("E2").Value = 12
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 11

$ws = $wb.Worksheets.Item("PEN")
$ws.Range("D2").Value = 11
$ws.Range("B3").Value = 48
$ws.Range("D3").Value = 9
$ws.Range("D4").Value = 14
